# Generate Report for Handoff
# Adds two new tracked files (two PNG images + the existing markdown file's
# handback refresh) to the localization-status workbook: one new row per
# file on the "Overview" sheet and on each language sheet ("zh-cn","de-de").

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276   # BGR packing of RGB FF6495ED (the workbook's HyperLink font color)

function Style-LinkCell($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $HYPERLINK_COLOR
}

function Set-Link($ws, $cellRef, $address, $display) {
    $rng = $ws.Range($cellRef)
    # NOTE: Range.Hyperlinks.Delete() wipes out *every* hyperlink on the
    # worksheet if called on a range that has none of its own - only call it
    # when this specific cell already owns a hyperlink to replace.
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $address, "", "", $display)
    Style-LinkCell $rng
}

# ---------------------------------------------------------------------
# Source data for the three tracked files
# ---------------------------------------------------------------------
$png1Name   = "3dfc73e4-932c-4f23-bbaa-d450118fdf6d.png"
$png1Target = "1f5ee9ff99f3938d4d6f1be35e4f1050230f1dc5.png"

$png2Name   = "dee78f61-f1a3-4412-9287-ca576f8ab42f.png"
$png2Target = "46a62874506a1fb33d41ba511e596ce4ddc9fb0a.png"

$mdName      = "edfe2ab8-e09e-4cda-8e62-9d65657fe69f.md"
$mdZhTarget  = "edfe2ab8-e09e-4cda-8e62-9d65657fe69f.402e9b1f57836e1368e3950740a3f1c7a422d473.zh-cn.xlf"
$mdDeTarget  = "edfe2ab8-e09e-4cda-8e62-9d65657fe69f.402e9b1f57836e1368e3950740a3f1c7a422d473.de-de.xlf"

$status       = "Ready for handoff"
$overviewDate = "2016-53-13 00:53:46"
$zhHandoffDt  = "2016-03-13 00:53:43"
$deHandoffDt  = "2016-03-13 00:53:46"
$epoch        = "0001-01-01 00:00:00"

$repoBase      = "https://github.com/OpenLocalizationTest/oltest/blob/b6ffd8db2195a6c5dd3f1b103aaa7814a72318cb/e2e"
$handoffZhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cfddc93d944da2099bff7d3e568145cd03ac6a11/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$handoffDeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c786026a29adbc5d60ea34d026d1447b97fa6a26/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# ===========================================================================
# Overview sheet
# ===========================================================================
$ov = $wb.Worksheets.Item("Overview")

# Refresh row 2 (previously the only tracked file, now the first png)
$ov.Range("A2").Value = $png1Name
$ov.Range("B2").Value = $status
$ov.Range("C2").Value = $status
$ov.Range("D2").Value = $overviewDate
Set-Link $ov "A2" ($repoBase + "/" + $png1Name) $png1Name

# Row 3 - second png
$ov.Range("A3").Value = $png2Name
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status
$ov.Range("D3").Value = $overviewDate
Set-Link $ov "A3" ($repoBase + "/" + $png2Name) $png2Name

# Row 4 - markdown file
$ov.Range("A4").Value = $mdName
$ov.Range("B4").Value = $status
$ov.Range("C4").Value = $status
$ov.Range("D4").Value = $overviewDate
Set-Link $ov "A4" ($repoBase + "/" + $mdName) $mdName

# ===========================================================================
# zh-cn sheet
# ===========================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 - first png
$zh.Range("A2").Value = $png1Name
$zh.Range("B2").Value = ".png"
$zh.Range("C2").Value = $status
$zh.Range("D2").Value = $png1Target
$zh.Range("E2").Value = $zhHandoffDt
$zh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H2").Value = $epoch
$zh.Range("I2").Value = "IsDependency"
$zh.Range("J2").Value = "e2e\" + $mdName
Set-Link $zh "A2" ($repoBase + "/" + $png1Name) $png1Name
Set-Link $zh "B2" ($repoBase + "/" + $png1Name) ".png"
Set-Link $zh "D2" ($handoffZhBase + "/" + $png1Target) $png1Target

# Row 3 - second png
$zh.Range("A3").Value = $png2Name
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $png2Target
$zh.Range("E3").Value = $zhHandoffDt
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = $epoch
$zh.Range("I3").Value = "IsDependency"
$zh.Range("J3").Value = "e2e\" + $mdName
Set-Link $zh "A3" ($repoBase + "/" + $png2Name) $png2Name
Set-Link $zh "B3" ($repoBase + "/" + $png2Name) ".png"
Set-Link $zh "D3" ($handoffZhBase + "/" + $png2Target) $png2Target

# Row 4 - markdown file
$zh.Range("A4").Value = $mdName
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = $status
$zh.Range("D4").Value = $mdZhTarget
$zh.Range("E4").Value = $zhHandoffDt
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = $epoch
$zh.Range("I4").Value = "Include"
Set-Link $zh "A4" ($repoBase + "/" + $mdName) $mdName
Set-Link $zh "B4" ($repoBase + "/" + $mdName) ".md"
Set-Link $zh "D4" ($handoffZhBase + "/" + $mdZhTarget) $mdZhTarget

# ===========================================================================
# de-de sheet
# ===========================================================================
$de = $wb.Worksheets.Item("de-de")

# Row 2 - first png
$de.Range("A2").Value = $png1Name
$de.Range("B2").Value = ".png"
$de.Range("C2").Value = $status
$de.Range("D2").Value = $png1Target
$de.Range("E2").Value = $deHandoffDt
$de.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H2").Value = $epoch
$de.Range("I2").Value = "IsDependency"
$de.Range("J2").Value = "e2e\" + $mdName
Set-Link $de "A2" ($repoBase + "/" + $png1Name) $png1Name
Set-Link $de "B2" ($repoBase + "/" + $png1Name) ".png"
Set-Link $de "D2" ($handoffDeBase + "/" + $png1Target) $png1Target

# Row 3 - second png
$de.Range("A3").Value = $png2Name
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = $status
$de.Range("D3").Value = $png2Target
$de.Range("E3").Value = $deHandoffDt
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = $epoch
$de.Range("I3").Value = "IsDependency"
$de.Range("J3").Value = "e2e\" + $mdName
Set-Link $de "A3" ($repoBase + "/" + $png2Name) $png2Name
Set-Link $de "B3" ($repoBase + "/" + $png2Name) ".png"
Set-Link $de "D3" ($handoffDeBase + "/" + $png2Target) $png2Target

# Row 4 - markdown file
$de.Range("A4").Value = $mdName
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = $status
$de.Range("D4").Value = $mdDeTarget
$de.Range("E4").Value = $deHandoffDt
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = $epoch
$de.Range("I4").Value = "Include"
Set-Link $de "A4" ($repoBase + "/" + $mdName) $mdName
Set-Link $de "B4" ($repoBase + "/" + $mdName) ".md"
Set-Link $de "D4" ($handoffDeBase + "/" + $mdDeTarget) $mdDeTarget

Write-Host "Localization status report refreshed for handoff."
